$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 177.57143
$ws.Cells.Item(12, 9).Value = 222.66667
$ws.Cells.Item(12, 10).Value = 143.75
$ws.Cells.Item(12, 11).Value = 222.66667
$ws.Cells.Item(12, 12).Value = 143.75
$ws.Cells.Item(12, 13).Value = -52.66667000000001
$ws.Cells.Item(12, 14).Value = -483.75
$ws.Cells.Item(43, 8).Value = 4022.1538
$ws.Cells.Item(43, 10).Value = 3612.8572
$ws.Cells.Item(43, 12).Value = 3612.8572
$ws.Cells.Item(43, 14).Value = -3750.8572
$ws.Cells.Item(70, 8).Value = 2333.3333
$ws.Cells.Item(70, 9).Value = 2333.3333
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 11).Value = 6999.999899999999
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 13).Value = -6729.999899999999
$ws.Cells.Item(70, 14).ClearContents()
$ws.Cells.Item(73, 8).Value = 2333.3333
$ws.Cells.Item(73, 9).Value = 2333.3333
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 11).Value = 6999.999899999999
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 13).Value = -6063.999899999999
$ws.Cells.Item(73, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 1811.1666
$ws.Cells.Item(132, 9).Value = 1451.4736
$ws.Cells.Item(132, 10).Value = 3178
$ws.Cells.Item(132, 11).Value = 4354.4208
$ws.Cells.Item(132, 12).Value = 9534
$ws.Cells.Item(132, 13).Value = -1824.4208
$ws.Cells.Item(132, 14).Value = -14594
$ws.Cells.Item(135, 8).Value = 2399.8
$ws.Cells.Item(135, 9).Value = 666.6667
$ws.Cells.Item(135, 11).Value = 6000.0003
$ws.Cells.Item(135, 13).Value = -3465.0003
$ws.Cells.Item(137, 8).Value = 1497
$ws.Cells.Item(137, 10).Value = 1400
$ws.Cells.Item(137, 12).Value = 4200
$ws.Cells.Item(137, 14).Value = -9300
$ws.Cells.Item(138, 8).Value = 4579.7
$ws.Cells.Item(138, 9).Value = 5324.25
$ws.Cells.Item(138, 10).Value = 4083.3333
$ws.Cells.Item(138, 11).Value = 15972.75
$ws.Cells.Item(138, 12).Value = 12249.9999
$ws.Cells.Item(138, 13).Value = -10832.75
$ws.Cells.Item(138, 14).Value = -22529.9999
$ws.Cells.Item(141, 8).Value = 911.6
$ws.Cells.Item(141, 9).Value = 924.1111
$ws.Cells.Item(141, 11).Value = 2772.3333
$ws.Cells.Item(141, 13).Value = 2407.6667

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 3335.4
$ws.Cells.Item(2, 9).Value = 1669.75
$ws.Cells.Item(2, 11).Value = 1669.75
$ws.Cells.Item(2, 13).Value = -1556.75
$ws.Cells.Item(45, 8).Value = 3078.6191
$ws.Cells.Item(45, 9).Value = 1315.8572
$ws.Cells.Item(45, 11).Value = 1315.8572
$ws.Cells.Item(45, 13).Value = -938.8571999999999
$ws.Cells.Item(61, 8).Value = 1533
$ws.Cells.Item(61, 9).Value = 799.5
$ws.Cells.Item(61, 11).Value = 799.5
$ws.Cells.Item(61, 13).Value = -587.5
$ws.Cells.Item(74, 8).Value = 849
$ws.Cells.Item(74, 9).Value = 849
$ws.Cells.Item(74, 11).Value = 849
$ws.Cells.Item(74, 13).Value = 25
$ws.Cells.Item(77, 8).Value = 849
$ws.Cells.Item(77, 9).Value = 849
$ws.Cells.Item(77, 11).Value = 4245
$ws.Cells.Item(77, 13).Value = 123
$ws.Cells.Item(102, 8).Value = 3000
$ws.Cells.Item(102, 9).Value = 3000
$ws.Cells.Item(102, 11).Value = 3000
$ws.Cells.Item(102, 13).Value = -1378
$ws.Cells.Item(116, 8).Value = 3335.4
$ws.Cells.Item(116, 9).Value = 1669.75
$ws.Cells.Item(116, 11).Value = 1669.75
$ws.Cells.Item(116, 13).Value = 624.25
$ws.Cells.Item(122, 8).Value = 4638
$ws.Cells.Item(122, 10).Value = 7777
$ws.Cells.Item(122, 12).Value = 23331
$ws.Cells.Item(122, 14).Value = -28231
$ws.Cells.Item(136, 8).Value = 1533
$ws.Cells.Item(136, 9).Value = 799.5
$ws.Cells.Item(136, 11).Value = 2398.5
$ws.Cells.Item(136, 13).Value = 151.5

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 3335.4
$ws.Cells.Item(3, 9).Value = 1669.75
$ws.Cells.Item(3, 11).Value = 1669.75
$ws.Cells.Item(3, 13).Value = -1555.75
$ws.Cells.Item(20, 8).Value = 2739.4443
$ws.Cells.Item(20, 10).Value = 3319
$ws.Cells.Item(20, 12).Value = 3319
$ws.Cells.Item(20, 14).Value = -3813
$ws.Cells.Item(132, 8).Value = 0
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 14).ClearContents()

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1500
$ws.Cells.Item(16, 9).Value = 1500
$ws.Cells.Item(16, 11).Value = 1500
$ws.Cells.Item(16, 13).Value = -1213
$ws.Cells.Item(31, 8).Value = 3387.75
$ws.Cells.Item(31, 9).Value = 2423.2
$ws.Cells.Item(31, 11).Value = 2423.2
$ws.Cells.Item(31, 13).Value = -2128.2
$ws.Cells.Item(34, 8).Value = 3387.75
$ws.Cells.Item(34, 9).Value = 2423.2
$ws.Cells.Item(34, 11).Value = 2423.2
$ws.Cells.Item(34, 13).Value = -2221.2
$ws.Cells.Item(58, 8).Value = 4500
$ws.Cells.Item(58, 9).Value = 5000
$ws.Cells.Item(58, 10).Value = 4000
$ws.Cells.Item(58, 11).Value = 5000
$ws.Cells.Item(58, 12).Value = 4000
$ws.Cells.Item(58, 13).Value = -4797
$ws.Cells.Item(58, 14).Value = -4406
$ws.Cells.Item(105, 8).Value = 8550.538
$ws.Cells.Item(105, 9).Value = 10565.7
$ws.Cells.Item(105, 11).Value = 10565.7
$ws.Cells.Item(105, 13).Value = -8818.700000000001
$ws.Cells.Item(113, 8).Value = 1500
$ws.Cells.Item(113, 9).Value = 1500
$ws.Cells.Item(113, 11).Value = 1500
$ws.Cells.Item(113, 13).Value = 670
$ws.Cells.Item(136, 8).Value = 4500
$ws.Cells.Item(136, 9).Value = 5000
$ws.Cells.Item(136, 10).Value = 4000
$ws.Cells.Item(136, 11).Value = 15000
$ws.Cells.Item(136, 12).Value = 12000
$ws.Cells.Item(136, 13).Value = -12450
$ws.Cells.Item(136, 14).Value = -17100

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 661.625
$ws.Cells.Item(34, 9).Value = 741.8570999999999
$ws.Cells.Item(34, 10).Value = 100
$ws.Cells.Item(34, 11).Value = 2225.5713
$ws.Cells.Item(34, 12).Value = 300
$ws.Cells.Item(34, 13).Value = -2141.5713
$ws.Cells.Item(34, 14).Value = -468
$ws.Cells.Item(103, 8).Value = 839.25
$ws.Cells.Item(103, 9).Value = 122.25
$ws.Cells.Item(103, 11).Value = 366.75
$ws.Cells.Item(103, 13).Value = 512.25
$ws.Cells.Item(113, 8).Value = 1095.7778
$ws.Cells.Item(113, 9).Value = 536.5
$ws.Cells.Item(113, 10).Value = 1255.5714
$ws.Cells.Item(113, 11).Value = 1609.5
$ws.Cells.Item(113, 12).Value = 3766.7142
$ws.Cells.Item(113, 13).Value = 560.5
$ws.Cells.Item(113, 14).Value = -8106.7142
$ws.Cells.Item(131, 8).Value = 914.1539
$ws.Cells.Item(131, 9).Value = 1099.3334
$ws.Cells.Item(131, 11).Value = 3298.0002
$ws.Cells.Item(131, 13).Value = 1741.9998
$ws.Cells.Item(137, 8).Value = 3589.5908
$ws.Cells.Item(137, 9).Value = 1627.5
$ws.Cells.Item(137, 10).Value = 3785.8
$ws.Cells.Item(137, 11).Value = 4882.5
$ws.Cells.Item(137, 12).Value = 11357.4
$ws.Cells.Item(137, 13).Value = 217.5
$ws.Cells.Item(137, 14).Value = -21557.4
$ws.Cells.Item(138, 8).Value = 2831.5
$ws.Cells.Item(138, 9).Value = 2663.3333
$ws.Cells.Item(138, 10).Value = 2999.6667
$ws.Cells.Item(138, 11).Value = 7989.999899999999
$ws.Cells.Item(138, 12).Value = 8999.000100000001
$ws.Cells.Item(138, 13).Value = -2849.999899999999
$ws.Cells.Item(138, 14).Value = -19279.0001
$ws.Cells.Item(139, 8).Value = 4300
$ws.Cells.Item(139, 9).Value = 6250
$ws.Cells.Item(139, 11).Value = 18750
$ws.Cells.Item(139, 13).Value = -13610

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 6008
$ws.Cells.Item(43, 9).Value = 6008
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 11).Value = 6008
$ws.Cells.Item(43, 12).Value = 0
$ws.Cells.Item(43, 13).Value = -5857
$ws.Cells.Item(43, 14).ClearContents()
$ws.Cells.Item(135, 8).Value = 80000
$ws.Cells.Item(135, 10).Value = 80000
$ws.Cells.Item(135, 12).Value = 80000
$ws.Cells.Item(135, 14).Value = -90140
$ws.Cells.Item(138, 8).Value = 109993
$ws.Cells.Item(138, 10).Value = 109993
$ws.Cells.Item(138, 12).Value = 109993
$ws.Cells.Item(138, 14).Value = -120273

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 14).ClearContents()
$ws.Cells.Item(27, 8).Value = 0
$ws.Cells.Item(27, 10).Value = 0
$ws.Cells.Item(27, 12).Value = 0
$ws.Cells.Item(27, 14).ClearContents()
$ws.Cells.Item(46, 8).Value = 4413.793
$ws.Cells.Item(46, 9).Value = 4000
$ws.Cells.Item(46, 10).Value = 4800
$ws.Cells.Item(46, 11).Value = 4000
$ws.Cells.Item(46, 12).Value = 4800
$ws.Cells.Item(46, 13).Value = -3812
$ws.Cells.Item(46, 14).Value = -5176
$ws.Cells.Item(55, 8).Value = 516.25
$ws.Cells.Item(55, 10).Value = 583.5714
$ws.Cells.Item(55, 12).Value = 583.5714
$ws.Cells.Item(55, 14).Value = -929.5714
$ws.Cells.Item(132, 8).Value = 1582.3334
$ws.Cells.Item(132, 9).Value = 1582.3334
$ws.Cells.Item(132, 11).Value = 4747.0002
$ws.Cells.Item(132, 13).Value = -2217.0002

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(58, 8).Value = 20094
$ws.Cells.Item(58, 10).Value = 20094
$ws.Cells.Item(58, 12).Value = 20094
$ws.Cells.Item(58, 14).Value = -20710
$ws.Cells.Item(62, 8).Value = 0
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 13).ClearContents()
$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 13).ClearContents()
$ws.Cells.Item(81, 8).Value = 3969.6
$ws.Cells.Item(81, 9).Value = 3855.111
$ws.Cells.Item(81, 11).Value = 7710.222
$ws.Cells.Item(81, 13).Value = -6649.222
$ws.Cells.Item(84, 8).Value = 3969.6
$ws.Cells.Item(84, 9).Value = 3855.111
$ws.Cells.Item(84, 11).Value = 38551.11
$ws.Cells.Item(84, 13).Value = -33247.11
$ws.Cells.Item(100, 8).Value = 4101379
$ws.Cells.Item(100, 9).Value = 5362234
$ws.Cells.Item(100, 11).Value = 10724468
$ws.Cells.Item(100, 13).Value = -10723927
$ws.Cells.Item(132, 8).Value = 2407.375
$ws.Cells.Item(132, 9).Value = 1959.8334
$ws.Cells.Item(132, 11).Value = 5879.5002
$ws.Cells.Item(132, 13).Value = -3349.5002
